$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    $cell.Formula = '="2013-06-13"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
